# Update odds data in row 2-22 (columns F..AO) to match the latest Betfair snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 5.5
$ws.Cells.Item(2, 7).Value = 7.2
$ws.Cells.Item(2, 8).Value = 1.77
$ws.Cells.Item(2, 9).Value = 1.94
$ws.Cells.Item(2, 10).Value = 3
$ws.Cells.Item(2, 12).Value = 1.56
$ws.Cells.Item(2, 13).Value = 1.15
$ws.Cells.Item(2, 14).Value = 2.24
$ws.Cells.Item(2, 15).Value = 1.64
$ws.Cells.Item(2, 16).Value = 1.4
$ws.Cells.Item(2, 17).Value = 2.94
$ws.Cells.Item(2, 19).Value = 6.8
$ws.Cells.Item(2, 20).Value = 2.52
$ws.Cells.Item(2, 21).Value = 1.51
$ws.Cells.Item(2, 22).Value = 2.08
$ws.Cells.Item(2, 23).Value = 1.15
$ws.Cells.Item(2, 24).Value = 7.6
$ws.Cells.Item(2, 25).Value = 5.5
$ws.Cells.Item(2, 26).Value = 9.4
$ws.Cells.Item(2, 27).Value = 23
$ws.Cells.Item(2, 28).Value = 15
$ws.Cells.Item(2, 30).Value = 12.5
$ws.Cells.Item(2, 31).Value = 32
$ws.Cells.Item(2, 32).Value = 150
$ws.Cells.Item(2, 33).Value = 32
$ws.Cells.Item(2, 34).Value = 38
$ws.Cells.Item(2, 35).Value = 200
$ws.Cells.Item(2, 41).Value = 28
$ws.Cells.Item(3, 7).Value = 1.59
$ws.Cells.Item(3, 8).Value = 6.6
$ws.Cells.Item(3, 16).Value = 2.16
$ws.Cells.Item(3, 17).Value = 1.77
$ws.Cells.Item(3, 18).Value = 1.41
$ws.Cells.Item(3, 19).Value = 3.05
$ws.Cells.Item(3, 20).Value = 1.88
$ws.Cells.Item(3, 23).Value = 2.72
$ws.Cells.Item(3, 27).Value = 240
$ws.Cells.Item(3, 39).Value = 150
$ws.Cells.Item(3, 40).Value = 8.800000000000001
$ws.Cells.Item(3, 41).Value = 160
$ws.Cells.Item(4, 8).Value = 1.55
$ws.Cells.Item(4, 9).Value = 1.56
$ws.Cells.Item(4, 13).Value = 1.04
$ws.Cells.Item(4, 17).Value = 1.66
$ws.Cells.Item(4, 18).Value = 1.57
$ws.Cells.Item(4, 21).Value = 2.18
$ws.Cells.Item(4, 22).Value = 2.78
$ws.Cells.Item(4, 26).Value = 9.6
$ws.Cells.Item(4, 33).Value = 24
$ws.Cells.Item(5, 6).Value = 1.98
$ws.Cells.Item(5, 7).Value = 2.12
$ws.Cells.Item(5, 8).Value = 3.65
$ws.Cells.Item(5, 9).Value = 4.1
$ws.Cells.Item(5, 10).Value = 3.75
$ws.Cells.Item(5, 14).Value = 4.7
$ws.Cells.Item(5, 16).Value = 2.3
$ws.Cells.Item(5, 17).Value = 1.64
$ws.Cells.Item(5, 18).Value = 1.51
$ws.Cells.Item(5, 19).Value = 2.62
$ws.Cells.Item(5, 21).Value = 2.38
$ws.Cells.Item(5, 23).Value = 1.9
$ws.Cells.Item(5, 24).Value = 21
$ws.Cells.Item(5, 25).Value = 18.5
$ws.Cells.Item(5, 27).Value = 70
$ws.Cells.Item(5, 31).Value = 100
$ws.Cells.Item(5, 33).Value = 11
$ws.Cells.Item(5, 35).Value = 120
$ws.Cells.Item(5, 38).Value = 30
$ws.Cells.Item(5, 39).Value = 330
$ws.Cells.Item(5, 41).Value = 32
$ws.Cells.Item(6, 10).Value = 4.5
$ws.Cells.Item(6, 16).Value = 2.68
$ws.Cells.Item(6, 17).Value = 1.49
$ws.Cells.Item(6, 18).Value = 1.7
$ws.Cells.Item(6, 24).Value = 29
$ws.Cells.Item(6, 38).Value = 27
$ws.Cells.Item(7, 15).Value = 1.12
$ws.Cells.Item(7, 19).Value = 1.96
$ws.Cells.Item(7, 21).Value = 2.56
$ws.Cells.Item(7, 28).Value = 46
$ws.Cells.Item(7, 32).Value = 85
$ws.Cells.Item(8, 6).Value = 2.88
$ws.Cells.Item(8, 9).Value = 2.68
$ws.Cells.Item(8, 15).Value = 1.38
$ws.Cells.Item(8, 16).Value = 1.81
$ws.Cells.Item(8, 19).Value = 4
$ws.Cells.Item(8, 25).Value = 10.5
$ws.Cells.Item(8, 27).Value = 110
$ws.Cells.Item(8, 34).Value = 21
$ws.Cells.Item(9, 9).Value = 4.5
$ws.Cells.Item(9, 29).Value = 23
$ws.Cells.Item(9, 34).Value = 970
$ws.Cells.Item(9, 37).Value = 500
$ws.Cells.Item(9, 38).Value = 500
$ws.Cells.Item(10, 12).Value = 1.32
$ws.Cells.Item(10, 21).Value = 1.84
$ws.Cells.Item(11, 11).Value = 5.9
$ws.Cells.Item(11, 18).Value = 1.67
$ws.Cells.Item(11, 19).Value = 2.24
$ws.Cells.Item(11, 32).Value = 23
$ws.Cells.Item(11, 36).Value = 44
$ws.Cells.Item(12, 6).Value = 1.65
$ws.Cells.Item(12, 7).Value = 1.75
$ws.Cells.Item(12, 8).Value = 5
$ws.Cells.Item(12, 9).Value = 5.7
$ws.Cells.Item(12, 10).Value = 4.2
$ws.Cells.Item(12, 11).Value = 4.9
$ws.Cells.Item(12, 16).Value = 2.22
$ws.Cells.Item(12, 17).Value = 1.66
$ws.Cells.Item(12, 22).Value = 1.21
$ws.Cells.Item(12, 23).Value = 2.32
$ws.Cells.Item(12, 25).Value = 110
$ws.Cells.Item(12, 30).Value = 500
$ws.Cells.Item(12, 40).Value = 8.800000000000001
$ws.Cells.Item(13, 12).Value = 1.3
$ws.Cells.Item(13, 17).Value = 1.71
$ws.Cells.Item(13, 19).Value = 2.78
$ws.Cells.Item(13, 24).Value = 18.5
$ws.Cells.Item(13, 29).Value = 9.199999999999999
$ws.Cells.Item(14, 20).Value = 1.98
$ws.Cells.Item(14, 21).Value = 1.83
$ws.Cells.Item(14, 23).Value = 3.2
$ws.Cells.Item(14, 34).Value = 55
$ws.Cells.Item(15, 7).Value = 2.58
$ws.Cells.Item(15, 8).Value = 3.1
$ws.Cells.Item(16, 7).Value = 2.94
$ws.Cells.Item(16, 9).Value = 2.56
$ws.Cells.Item(16, 21).Value = 2.42
$ws.Cells.Item(16, 24).Value = 18
$ws.Cells.Item(16, 27).Value = 36
$ws.Cells.Item(16, 40).Value = 22
$ws.Cells.Item(17, 6).Value = 4.2
$ws.Cells.Item(17, 7).Value = 4.3
$ws.Cells.Item(17, 9).Value = 2.08
$ws.Cells.Item(17, 19).Value = 3.95
$ws.Cells.Item(17, 23).Value = 1.3
$ws.Cells.Item(17, 24).Value = 13.5
$ws.Cells.Item(17, 26).Value = 11.5
$ws.Cells.Item(17, 31).Value = 22
$ws.Cells.Item(17, 33).Value = 17
$ws.Cells.Item(17, 41).Value = 16.5
$ws.Cells.Item(18, 7).Value = 3.05
$ws.Cells.Item(18, 8).Value = 2.66
$ws.Cells.Item(18, 9).Value = 2.68
$ws.Cells.Item(18, 10).Value = 3.35
$ws.Cells.Item(18, 11).Value = 3.4
$ws.Cells.Item(18, 14).Value = 3.3
$ws.Cells.Item(18, 17).Value = 2.26
$ws.Cells.Item(18, 20).Value = 1.88
$ws.Cells.Item(18, 22).Value = 1.59
$ws.Cells.Item(18, 24).Value = 12
$ws.Cells.Item(18, 26).Value = 16
$ws.Cells.Item(18, 41).Value = 30
$ws.Cells.Item(19, 6).Value = 4.6
$ws.Cells.Item(19, 7).Value = 4.7
$ws.Cells.Item(19, 8).Value = 1.8
$ws.Cells.Item(19, 9).Value = 1.81
$ws.Cells.Item(19, 22).Value = 2.22
$ws.Cells.Item(19, 23).Value = 1.27
$ws.Cells.Item(19, 32).Value = 38
$ws.Cells.Item(20, 8).Value = 7
$ws.Cells.Item(20, 11).Value = 5.4
$ws.Cells.Item(20, 18).Value = 1.86
$ws.Cells.Item(20, 19).Value = 2.12
$ws.Cells.Item(20, 25).Value = 38
$ws.Cells.Item(20, 36).Value = 14.5
$ws.Cells.Item(20, 38).Value = 23
$ws.Cells.Item(21, 6).Value = 1.33
$ws.Cells.Item(21, 7).Value = 1.34
$ws.Cells.Item(21, 9).Value = 11
$ws.Cells.Item(21, 10).Value = 6.4
$ws.Cells.Item(21, 11).Value = 6.6
$ws.Cells.Item(21, 23).Value = 3.95
$ws.Cells.Item(21, 35).Value = 85
$ws.Cells.Item(21, 36).Value = 12.5
$ws.Cells.Item(22, 6).Value = 2.6
$ws.Cells.Item(22, 7).Value = 2.62
$ws.Cells.Item(22, 17).Value = 1.79
$ws.Cells.Item(22, 23).Value = 1.61
$ws.Cells.Item(22, 40).Value = 18
